$d = $word.ActiveDocument
$t = $d.Tables(1)

# --- Step 1: delete the four now-obsolete rows (original rows 8,9,10,11:
#     0.00000 / 0.00003 / 0.00003 / 0.00003). Delete highest index first
#     so earlier indices stay stable while deleting.
$t.Rows(11).Delete()
$t.Rows(10).Delete()
$t.Rows(9).Delete()
$t.Rows(8).Delete()

# --- Step 2: insert four new (blank) rows before what is now row 5
#     (still the original "0.00003" row, since rows 8-11 were after it).
$t.Rows.Add($t.Rows(5)) | Out-Null
$t.Rows.Add($t.Rows(5)) | Out-Null
$t.Rows.Add($t.Rows(5)) | Out-Null
$t.Rows.Add($t.Rows(5)) | Out-Null

# --- Step 3: set the text content for every row that changed.
$t.Rows(1).Cells(1).Range.Text = "0M"
$t.Rows(2).Cells(1).Range.Text = "0M"
$t.Rows(3).Cells(1).Range.Text = "0M"
$t.Rows(4).Cells(1).Range.Text = "63"

$t.Rows(5).Cells(1).Range.Text = "0.00002"
$t.Rows(6).Cells(1).Range.Text = "0.00007"
$t.Rows(7).Cells(1).Range.Text = "0.00004"
$t.Rows(8).Cells(1).Range.Text = "0.00001"

$t.Rows(9).Cells(1).Range.Text = "0.00004"

$t.Rows(12).Cells(1).Range.Text = "0.00219"

$t.Rows(44).Cells(1).Range.Text = "100"
$t.Rows(45).Cells(1).Range.Text = "0"
$t.Rows(46).Cells(1).Range.Text = "410"
